# Generate Report for Handoff
# Adds two new localization-status rows (one .md source, and two .png sources)
# to the Overview / zh-cn / de-de sheets, mirroring the existing row for the
# first file, and refreshes the "Latest Handoff Date" for the already
# present row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# New source files being reported on in this handoff batch
# ---------------------------------------------------------------------
$fileA = "7bf0f3da-3740-43d2-8cdb-1c39c885ed29.png"
$fileB = "d8063296-b9d0-4678-bb1c-ccaa65b2e1c7.md"
$fileC = "e9858d9e-af34-467c-aedc-694c44630512.png"

$srcRepo = "https://github.com/OpenLocalizationTest/oltest/blob/fbdb0c1deab485612d13ee47c27a7bc50ec714a5/e2e"
$urlA = "$srcRepo/$fileA"
$urlB = "$srcRepo/$fileB"
$urlC = "$srcRepo/$fileC"

$handoffDate = "2016-49-18 22:49:06"
$handoffDatetime = "2016-03-18 22:49:02"
$epoch = "0001-01-01 00:00:00"
$status = "Ready for handoff"

# zh-cn handoff targets
$zhRepo = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5c1e8f2a93b7d4016e2fa0c8b3d7a45e9f2c6d81/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
$pngTargetA_zh = "d7076cbb31f225024d830929a5261d1cd53c8f59.png"
$xlfTargetB_zh = "d8063296-b9d0-4678-bb1c-ccaa65b2e1c7.4c3d118645c9ee5123274310aea455455f01f0ad.zh-cn.xlf"
$pngTargetC_zh = "8a36f39968d27098ae52df1c0a7225ff3fc89b22.png"
$zhDate = "2016-03-18 22:49:02"

# de-de handoff targets
$deRepo = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9a3d6c1f84e2b5709a1c4e7b2d6f9a3c8e5b1d07/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"
$pngTargetA_de = "d7076cbb31f225024d830929a5261d1cd53c8f59.png"
$xlfTargetB_de = "d8063296-b9d0-4678-bb1c-ccaa65b2e1c7.4c3d118645c9ee5123274310aea455455f01f0ad.de-de.xlf"
$pngTargetC_de = "8a36f39968d27098ae52df1c0a7225ff3fc89b22.png"
$deDate = "2016-03-18 22:49:06"

# ======================================================================
# Sheet "Overview"
# ======================================================================
$ov = $wb.Worksheets.Item("Overview")

# Refresh the handoff date for the already-reported file (row 2)
$ov.Range("D2").Value = $handoffDate

# Row 3: the new .md file
$ov.Range("B3").Value = $status
$ov.Range("C3").Value = $status
$ov.Range("D3").Value = $handoffDate

# Row 4: the new .png file
$ov.Range("B4").Value = $status
$ov.Range("C4").Value = $status
$ov.Range("D4").Value = $handoffDate

# Rebuild every hyperlink on the sheet (row 2's target file was renamed from
# .md to .png upstream, so its display text/url need refreshing too) in
# row order so relationship ids come out sequential.
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), $urlA, "", "", $fileA) | Out-Null
$ov.Hyperlinks.Add($ov.Range("A3"), $urlB, "", "", $fileB) | Out-Null
$ov.Hyperlinks.Add($ov.Range("A4"), $urlC, "", "", $fileC) | Out-Null

# ======================================================================
# Sheet "zh-cn"
# ======================================================================
$zh = $wb.Worksheets.Item("zh-cn")

# --- Row 2 (existing file, now reported as a .png) ---
$zh.Range("B2").Value = ".png"
$zh.Range("D2").Value = $pngTargetA_zh
$zh.Range("E2").Value = $zhDate
$zh.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("H2").Value = $epoch
$zh.Range("I2").Value = "IsDependency"
$zh.Range("J2").Value = "e2e\$fileB"

# --- Row 3 (new .md file) ---
$zh.Range("C3").Value = $status
$zh.Range("E3").Value = $zhDate
$zh.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("H3").Value = $epoch
$zh.Range("I3").Value = "Include"

# --- Row 4 (new .png file) ---
$zh.Range("C4").Value = $status
$zh.Range("E4").Value = $zhDate
$zh.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("H4").Value = $epoch
$zh.Range("I4").Value = "IsDependency"
$zh.Range("J4").Value = "e2e\$fileB"

# Rebuild every hyperlink on the sheet, row by row, columns A / B / D
# (matches the relationship-id allocation order in the target workbook).
$zh.Hyperlinks.Delete()

$zh.Hyperlinks.Add($zh.Range("A2"), $urlA, "", "", $fileA) | Out-Null
$zh.Hyperlinks.Add($zh.Range("B2"), $urlA, "", "", ".png") | Out-Null
$zh.Hyperlinks.Add($zh.Range("D2"), "$zhRepo/$pngTargetA_zh", "", "", $pngTargetA_zh) | Out-Null

$zh.Hyperlinks.Add($zh.Range("A3"), $urlB, "", "", $fileB) | Out-Null
$zh.Hyperlinks.Add($zh.Range("B3"), $urlB, "", "", ".md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("D3"), "$zhRepo/$xlfTargetB_zh", "", "", $xlfTargetB_zh) | Out-Null

$zh.Hyperlinks.Add($zh.Range("A4"), $urlC, "", "", $fileC) | Out-Null
$zh.Hyperlinks.Add($zh.Range("B4"), $urlC, "", "", ".png") | Out-Null
$zh.Hyperlinks.Add($zh.Range("D4"), "$zhRepo/$pngTargetC_zh", "", "", $pngTargetC_zh) | Out-Null

# ======================================================================
# Sheet "de-de"
# ======================================================================
$de = $wb.Worksheets.Item("de-de")

# --- Row 2 (existing file, now reported as a .png) ---
$de.Range("B2").Value = ".png"
$de.Range("D2").Value = $pngTargetA_de
$de.Range("E2").Value = $deDate
$de.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("H2").Value = $epoch
$de.Range("I2").Value = "IsDependency"
$de.Range("J2").Value = "e2e\$fileB"

# --- Row 3 (new .md file) ---
$de.Range("C3").Value = $status
$de.Range("E3").Value = $deDate
$de.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("H3").Value = $epoch
$de.Range("I3").Value = "Include"

# --- Row 4 (new .png file) ---
$de.Range("C4").Value = $status
$de.Range("E4").Value = $deDate
$de.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("H4").Value = $epoch
$de.Range("I4").Value = "IsDependency"
$de.Range("J4").Value = "e2e\$fileB"

# Rebuild every hyperlink on the sheet, row by row, columns A / B / D
# (matches the relationship-id allocation order in the target workbook).
$de.Hyperlinks.Delete()

$de.Hyperlinks.Add($de.Range("A2"), $urlA, "", "", $fileA) | Out-Null
$de.Hyperlinks.Add($de.Range("B2"), $urlA, "", "", ".png") | Out-Null
$de.Hyperlinks.Add($de.Range("D2"), "$deRepo/$pngTargetA_de", "", "", $pngTargetA_de) | Out-Null

$de.Hyperlinks.Add($de.Range("A3"), $urlB, "", "", $fileB) | Out-Null
$de.Hyperlinks.Add($de.Range("B3"), $urlB, "", "", ".md") | Out-Null
$de.Hyperlinks.Add($de.Range("D3"), "$deRepo/$xlfTargetB_de", "", "", $xlfTargetB_de) | Out-Null

$de.Hyperlinks.Add($de.Range("A4"), $urlC, "", "", $fileC) | Out-Null
$de.Hyperlinks.Add($de.Range("B4"), $urlC, "", "", ".png") | Out-Null
$de.Hyperlinks.Add($de.Range("D4"), "$deRepo/$pngTargetC_de", "", "", $pngTargetC_de) | Out-Null

Write-Host "Handoff report rows added."
